# [Kadastro App] Yeni kayit eklendi: 2923
# Adds a new record (Kayit No 2923) to the "Kayitlar" master sheet and to the
# filtered "Erdemli" district sheet (the record's Birim is "Erdemli").

$wb = $excel.ActiveWorkbook

function Add-KayitRow($SheetName, $RowIndex) {
    $ws = $wb.Worksheets.Item($SheetName)

    $rng = $ws.Range("A" + $RowIndex + ":F" + $RowIndex)

    # Use a quote-prefix so values that look numeric/date-like ("2923", "2025-09-08",
    # "2") are kept as plain text, matching the rest of the sheet (all cells are
    # stored as text strings).
    $ws.Range("A" + $RowIndex).Value = "'2923"
    $ws.Range("B" + $RowIndex).Value = "'2025-09-08"
    $ws.Range("C" + $RowIndex).Value = "Erdemli"
    $ws.Range("D" + $RowIndex).Value = "'2"
    $ws.Range("E" + $RowIndex).Value = "LİHKAB"
    $ws.Range("F" + $RowIndex).Value = "EMİNE ALANLI KIRCILI (K.Mühendisi), CEMAL TİMUROĞLU (K.Teknisyeni)"

    # Drop the quote-prefix/text formatting marker so the new cells end up with
    # the same (default) style as every other cell in the sheet.
    $rng.Style = "Normal"
}

# "Kayitlar" is the master list of all records -> new row 17
Add-KayitRow "Kayitlar" 17

# "Erdemli" is filtered to Birim = Erdemli -> new row 16
Add-KayitRow "Erdemli" 16
